$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "DT"
$ws.Range("C2").Value = 0.73568281938326
$ws.Range("D2").Value = 0.7354836205558143
$ws.Range("E2").Value = 0.7351172480739726
$ws.Range("F2").Value = 0.7354835935802485
$ws.Range("G2").Value = 0.73568281938326
$ws.Range("H2").Value = 0.7352844217283686
$ws.Range("I2").Value = 0.2647155782716314
$ws.Range("J2").Value = 0.26431718061674

$ws.Range("B3").Value = "KNN"
$ws.Range("C3").Value = 0.7004405286343612
$ws.Range("D3").Value = 0.700279459711225
$ws.Range("E3").Value = 0.7000213037482952
$ws.Range("F3").Value = 0.7002794411877641
$ws.Range("G3").Value = 0.7004405286343612
$ws.Range("H3").Value = 0.7001183907880888
$ws.Range("I3").Value = 0.2998816092119112
$ws.Range("J3").Value = 0.2995594713656388

$ws.Range("B4").Value = "GNB"
$ws.Range("C4").Value = 0.7929515418502202
$ws.Range("D4").Value = 0.792889302903276
$ws.Range("E4").Value = 0.7929033097243919
$ws.Range("F4").Value = 0.7928893004605095
$ws.Range("G4").Value = 0.7929515418502202
$ws.Range("H4").Value = 0.7928270639563314
$ws.Range("I4").Value = 0.2071729360436686
$ws.Range("J4").Value = 0.2070484581497798

$ws.Range("B5").Value = "SVM"
$ws.Range("C5").Value = 0.6916299559471366
$ws.Range("D5").Value = 0.691585157584226
$ws.Range("E5").Value = 0.6915940432867266
$ws.Range("F5").Value = 0.6915851561332889
$ws.Range("G5").Value = 0.6916299559471366
$ws.Range("H5").Value = 0.6915403592213155
$ws.Range("I5").Value = 0.3084596407786845
$ws.Range("J5").Value = 0.3083700440528634

$ws.Range("B6").Value = "LR"
$ws.Range("C6").Value = 0.7929515418502202
$ws.Range("D6").Value = 0.7928504890544946
$ws.Range("E6").Value = 0.792830914701977
$ws.Range("F6").Value = 0.7928504826146502
$ws.Range("G6").Value = 0.7929515418502202
$ws.Range("H6").Value = 0.792749436258769
$ws.Range("I6").Value = 0.207250563741231
$ws.Range("J6").Value = 0.2070484581497798

$ws.Range("B7").Value = "MLP"
$ws.Range("C7").Value = 0.5022026431718062
$ws.Range("D7").Value = 0.5
$ws.Range("E7").Value = 0.335783585463847
$ws.Range("F7").Value = 0.4999951483395191
$ws.Range("G7").Value = 0.5022026431718062
$ws.Range("H7").Value = 0.4977973568281938
$ws.Range("I7").Value = 0.5022026431718062
$ws.Range("J7").Value = 0.4977973568281938

$ws.Range("B8").Value = "XGB"
$ws.Range("C8").Value = 0.748898678414097
$ws.Range("D8").Value = 0.7487191429902189
$ws.Range("E8").Value = 0.7484592547631853
$ws.Range("F8").Value = 0.7487191214648117
$ws.Range("G8").Value = 0.748898678414097
$ws.Range("H8").Value = 0.7485396075663409
$ws.Range("I8").Value = 0.2514603924336591
$ws.Range("J8").Value = 0.251101321585903

$ws.Range("B9").Value = "RF"
$ws.Range("C9").Value = 0.73568281938326
$ws.Range("D9").Value = 0.7352895513119081
$ws.Range("E9").Value = 0.7335098391166802
$ws.Range("F9").Value = 0.7352894461425998
$ws.Range("G9").Value = 0.73568281938326
$ws.Range("H9").Value = 0.7348962832405563
$ws.Range("I9").Value = 0.2651037167594437
$ws.Range("J9").Value = 0.26431718061674

$ws.Range("B10").Value = "ET"
$ws.Range("C10").Value = 0.7533039647577092
$ws.Range("D10").Value = 0.7531439217512809
$ws.Range("E10").Value = 0.7529587207338903
$ws.Range("F10").Value = 0.7531439047467194
$ws.Range("G10").Value = 0.7533039647577092
$ws.Range("H10").Value = 0.7529838787448525
$ws.Range("I10").Value = 0.2470161212551475
$ws.Range("J10").Value = 0.2466960352422908

$ws.Range("B11").Value = "AdaBoost"
$ws.Range("C11").Value = 0.7224669603524229
$ws.Range("D11").Value = 0.7222869119701909
$ws.Range("E11").Value = 0.7219812815803627
$ws.Range("F11").Value = 0.7222868895293715
$ws.Range("G11").Value = 0.7224669603524229
$ws.Range("H11").Value = 0.7221068635879591
$ws.Range("I11").Value = 0.2778931364120409
$ws.Range("J11").Value = 0.2775330396475771

$ws.Range("B12").Value = "GB"
$ws.Range("C12").Value = 0.7224669603524229
$ws.Range("D12").Value = 0.722442167365316
$ws.Range("E12").Value = 0.7224561875997922
$ws.Range("F12").Value = 0.7224421669398894
$ws.Range("G12").Value = 0.7224669603524229
$ws.Range("H12").Value = 0.7224173743782091
$ws.Range("I12").Value = 0.2775826256217909
$ws.Range("J12").Value = 0.2775330396475771

$ws.Range("B13").Value = "SGD"
$ws.Range("C13").Value = 0.4933920704845815
$ws.Range("D13").Value = 0.493828598043782
$ws.Range("E13").Value = 0.4883701562454827
$ws.Range("F13").Value = 0.4938284051060422
$ws.Range("G13").Value = 0.4933920704845815
$ws.Range("H13").Value = 0.4942651256029826
$ws.Range("I13").Value = 0.5057348743970174
$ws.Range("J13").Value = 0.5066079295154184

$ws.Rows.Item(14).Delete()
